$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.09280751105506574
$ws.Range("B3").Value = 0.103151799219635
$ws.Range("H3").Value = 0.1959593102747007
$ws.Range("B4").Value = 0.09048863435816971
$ws.Range("H4").Value = 0.1832961454132355
$ws.Range("B5").Value = 0.06387097079667473
$ws.Range("C5").Value = 0.006256020468073837
$ws.Range("D5").Value = 10.61619576591366
$ws.Range("E5").Value = 0.06818126676673325
$ws.Range("F5").Value = 0.05158639374718649
$ws.Range("G5").Value = 0.07615554784616198
$ws.Range("H5").Value = 0.1566784818517405
$ws.Range("B6").Value = 0.04005464743891991
$ws.Range("C6").Value = 0.004456084037186648
$ws.Range("D6").Value = 4.713788261422956
$ws.Range("E6").Value = 0.03321303009358303
$ws.Range("F6").Value = 0.03131215345572928
$ws.Range("G6").Value = 0.04879714142210963
$ws.Range("H6").Value = 0.1328621584939856
$ws.Range("B7").Value = 0.03111834027607828
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = 0.123925851331144
$ws.Range("B8").Value = 0.02732859978610495
$ws.Range("C8").Value = 0.002703184104130505
$ws.Range("D8").Value = 2.616546450868733
$ws.Range("E8").Value = 0.01354470906404969
$ws.Range("F8").Value = 0.02202435973850002
$ws.Range("G8").Value = 0.0326328398337099
$ws.Range("H8").Value = 0.1201361108411707
$ws.Range("B9").Value = 0.02521493548482489
$ws.Range("C9").Value = 0.002887811832640039
$ws.Range("D9").Value = 2.13909171682071
$ws.Range("E9").Value = 0.01819187325317444
$ws.Range("F9").Value = 0.01954940697130819
$ws.Range("G9").Value = 0.0308804639983412
$ws.Range("H9").Value = 0.1180224465398906
$ws.Range("B10").Value = 0.02450256121722537
$ws.Range("C10").Value = 0.002623639709169693
$ws.Range("D10").Value = 1.884502709513444
$ws.Range("E10").Value = 0.01780690246746566
$ws.Range("F10").Value = 0.01935504060988099
$ws.Range("G10").Value = 0.02965008182456945
$ws.Range("H10").Value = 0.1173100722722911
$ws.Range("B11").Value = 0.03298186556308941
$ws.Range("H11").Value = 0.1257893766181551
$ws.Range("B12").Value = 0.0491487951751833
$ws.Range("H12").Value = 0.141956306230249
$ws.Range("B13").Value = 0.06100490878150903
$ws.Range("H13").Value = 0.1538124198365748
$ws.Range("B14").Value = 0.06606352031034302
$ws.Range("H14").Value = 0.1588710313654088
$ws.Range("B15").Value = 0.07210076590075497
$ws.Range("H15").Value = 0.1649082769558207
$ws.Range("B16").Value = 0.07739935413739364
$ws.Range("H16").Value = 0.1702068651924594
$ws.Range("B17").Value = 0.07893862506383548
$ws.Range("H17").Value = 0.1717461361189012
$ws.Range("B18").Value = -0.09280751105506574
$ws.Range("C18").Value = 0.01092901732207689
$ws.Range("D18").Value = -16.51209139823487
$ws.Range("E18").Value = 0.04191690866172165
$ws.Range("F18").Value = -0.1142767561183235
$ws.Range("G18").Value = -0.07133826599180827
$ws.Range("B19").Value = 0.07884046889105732
$ws.Range("H19").Value = 0.1716479799461231
$ws.Range("B20").Value = 0.08467423743518479
$ws.Range("H20").Value = 0.1774817484902505
$ws.Range("B21").Value = 0.08772785525910158
$ws.Range("H21").Value = 0.1805353663141673
$ws.Range("B22").Value = 0.09058560586556964
$ws.Range("H22").Value = 0.1833931169206354
$ws.Range("B23").Value = 0.09369842156729383
$ws.Range("H23").Value = 0.1865059326223596
$ws.Range("B24").Value = 0.09608927244233402
$ws.Range("C24").Value = 0.008627719292247897
$ws.Range("D24").Value = 223035066761.7094
$ws.Range("E24").Value = 0.05386315405091487
$ws.Range("F24").Value = 0.07913560427109874
$ws.Range("G24").Value = 0.1130429406135694
$ws.Range("H24").Value = 0.1888967834973997
$ws.Range("B25").Value = 0.09713797541985594
$ws.Range("C25").Value = 0.008427197005246143
$ws.Range("D25").Value = 465908856722.3784
$ws.Range("E25").Value = 0.05166444712238601
$ws.Range("F25").Value = 0.08057705399739304
$ws.Range("G25").Value = 0.1136988968423189
$ws.Range("H25").Value = 0.1899454864749217
$ws.Range("B26").Value = 0.100270206048319
$ws.Range("C26").Value = 0.008697650599882477
$ws.Range("D26").Value = 590819431733.5083
$ws.Range("E26").Value = 0.05710141576481517
$ws.Range("F26").Value = 0.08318171286105486
$ws.Range("G26").Value = 0.1173586992355835
$ws.Range("H26").Value = 0.1930777171033847
$ws.Range("B27").Value = 0.1053440025506479
$ws.Range("C27").Value = 0.008939364405148918
$ws.Range("D27").Value = 587199074062.0693
$ws.Range("E27").Value = 0.04922793667344875
$ws.Range("F27").Value = 0.08777675243249138
$ws.Range("G27").Value = 0.1229112526688048
$ws.Range("H27").Value = 0.1981515136057136
$ws.Range("B28").Value = 0.1017496236079487
$ws.Range("C28").Value = 0.008846132272907111
$ws.Range("D28").Value = 16.87761100610693
$ws.Range("E28").Value = 0.07632351879072788
$ws.Range("F28").Value = 0.08437087950573437
$ws.Range("G28").Value = 0.1191283677101631
$ws.Range("H28").Value = 0.1945571346630145
$ws.Range("B29").Value = 0.0304729534264252
$ws.Range("C29").Value = 0.003074019037222985
$ws.Range("D29").Value = 3.384110462072563
$ws.Range("E29").Value = 0.01694674548598612
$ws.Range("F29").Value = 0.02443253681086965
$ws.Range("G29").Value = 0.03651337004198098
$ws.Range("H29").Value = 0.1232804644814909
